$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.389.63'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.370.43'
$ws.Range("E3").Value = '  +5.35%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '232.92'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.19'
$ws.Range("E7").Value = '  +7.71%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.460'
$ws.Range("E9").Value = '  +2.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0955'
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.87'
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.65'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.723.28'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.73'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.26'
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.376.91'
$ws.Range("E18").Value = '  +5.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.420.73'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0989'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.22'
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.28'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.93'
$ws.Range("E24").Value = '  +16.63%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.46'
$ws.Range("E29").Value = '  +7.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.08'
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.51'
$ws.Range("E31").Value = '  +9.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("E32").Value = '  -6.15%  '
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("E34").Value = '  +4.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0691'
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("E36").Value = '  +2.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.51'
$ws.Range("E37").Value = '  +9.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.53'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0256'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.02'
$ws.Range("E41").Value = '  +9.63%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.22'
$ws.Range("E43").Value = '  +5.08%  '
$ws.Range("E44").Value = '  +8.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.38'
$ws.Range("E45").Value = '  +2.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0952'
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.447.73'
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.596.28'
$ws.Range("E50").Value = '  +5.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000204'
$ws.Range("E51").Value = '  -8.35%  '
